$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "myMvc790"
$ws.Range("B2").Value = 23091356
$ws.Range("C2").Value = "rijbqfl34"
$ws.Range("D2").Value = 'S8n$g2P&'
$ws.Range("F2").Value = "LFTclhtR"
$ws.Range("G2").Value = "kNDg"
